$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# Add new shared string "capacity" - handled automatically by engine when cell value assigned.

# Update header row (row 1)
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Update data row (row 2)
$ws.Range("B2").Value = "曰產"
$ws.Range("C2").Value = 1998
$ws.Range("D2").Value = "周桂香"
$ws.Range("E2").Value = "92年10月15曰"
$ws.Range("F2").Value = "買賣"
$ws.Range("G2").Value = "(超過五年）"
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2012-04-27"
$ws.Range("K2").Value = "許智傑"
$ws.Range("L2").Value = 1750
$ws.Range("M2").Value = "tmp1dd71"
$ws.Range("N2").Value = 38
